$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("potential_preg_untrt")
$ws.Range("C9").Value  = 0.05
$ws.Range("C10").Value = 0.02
$ws.Range("C11").Value = 0.02
$ws.Range("C13").Value = 0.005
$ws.Range("C14").Value = 0.004
$ws.Range("C15").Value = 0.004
$ws.Range("C16").Value = 0.004
$ws.Range("C17").Value = 0.004

# Active sheet / tab selection changes: SimParameters loses tabSelected,
# potential_preg_untrt becomes the active/selected sheet (3rd tab, index 2).
$ws.Activate()
$ws.Range("I16").Select()
